$d = $word.ActiveDocument
Write-Output "bookmarks count: $($d.Bookmarks.Count)"
foreach ($bm in $d.Bookmarks) {
    Write-Output "bm: $($bm.Name) start=$($bm.Start) end=$($bm.End)"
}
